$d = $word.ActiveDocument

$pairs = @(
    @("2025-03-10 Monday", "2025-03-11 Tuesday"),
    @("815×4=3260", "647×7=4529"),
    @("174×9=1566", "443×6=2658"),
    @("413×8=3304", "116×2=232"),
    @("625×9=5625", "137×7=959"),
    @("706×9=6354", "799×4=3196"),
    @("571×5=2855", "776×6=4656"),
    @("521×2=1042", "191×7=1337"),
    @("739×2=1478", "621×2=1242"),
    @("741×9=6669", "352×8=2816"),
    @("880×3=2640", "731×6=4386"),
    @("285×3=855", "336×9=3024"),
    @("835×8=6680", "177×9=1593"),
    @("211×8=1688", "257×9=2313"),
    @("327×5=1635", "535×6=3210"),
    @("597×7=4179", "674×4=2696"),
    @("599×7=4193", "976×5=4880"),
    @("159×5=795", "303×2=606"),
    @("531×6=3186", "213×8=1704"),
    @("413×7=2891", "871×6=5226"),
    @("743×5=3715", "544×4=2176"),
    @("963×4=3852", "855×4=3420"),
    @("673×7=4711", "303×9=2727"),
    @("904×2=1808", "857×9=7713"),
    @("133×8=1064", "175×5=875"),
    @("735×7=5145", "429×7=3003")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
